$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data rows 2 and 3 for columns D, J, K, L, M, P
# (columns A, B, C, E, F, G, H, I, N, O, Q, R are identical between the two rows)

$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $addr2 = "$col`2"
    $addr3 = "$col`3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}
